$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Insert a new "Meta description" paragraph right after the
# opening Heading1 paragraph ("Play Book of Shadows for Free - Review of
# Gameplay & Bonus Features").
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item(1)
$insertionPoint = $d.Range($headingPara.Range.End, $headingPara.Range.End)
$insertionPoint.InsertParagraphBefore()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"
$metaParaStart = $metaPara.Range.Start

$metaLabel = "Meta description"
$metaText = ": Read our review of Book of Shadows and play for free today. Learn about the gameplay mechanics and bonus features, including the Free Spins and Nolimit bonus."
$metaPara.Range.Text = $metaLabel + $metaText

# Bold just the "Meta description" label, leaving the rest regular.
$boldRange = $d.Range($metaParaStart, $metaParaStart + $metaLabel.Length)
$boldRange.Bold = 1

# ---------------------------------------------------------------------------
# Change 2: Remove the duplicate bold "Play Book of Shadows..." paragraph
# that used to sit just before the closing italic paragraph, and replace the
# text of that trailing italic paragraph with the new image-prompt text
# (keeping its italic formatting intact).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupHeadingPara = $d.Paragraphs.Item($count - 1)
$dupHeadingPara.Range.Delete()

$count = $d.Paragraphs.Count
$imagePromptPara = $d.Paragraphs.Item($count)
$imagePromptRange = $d.Range($imagePromptPara.Range.Start, $imagePromptPara.Range.End - 1)
$imagePromptRange.Text = "Create a feature image for Book of Shadows, a horror-themed slot game, that is in cartoon style and features a happy Maya warrior with glasses. The image should showcase the Maya warrior holding the Book of Shadows with a confident and mischievous expression on his face, ready to tackle the horrors in the game. The background should depict a dark, eerie forest with moonlight casting a shadowy glow. The Maya warrior should be depicted wearing traditional Maya clothing, including a headdress adorned with feathers. The glasses should be modern, adding a playful touch to the image. The overall style should be a mix of ancient and modern, representing the theme of the game. The image should be colorful, bold, and attention-grabbing to entice players to try the game."

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
